$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing city values (row 2-4)
$ws.Range("A2").Value = "London"
$ws.Range("A3").Value = "Amsterdan"
$ws.Range("A4").Value = "Paris"

# Add two new cities (row 5-6)
$ws.Range("A5").Value = "Rome"
$ws.Range("A6").Value = "Madrid"

# Update the selection to match the target state
$ws.Range("A2").Select()
